$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.017.63"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.760.49"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'575.96"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").Value = "'159.16"
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -3.43%  "
$ws.Range("E9").Value = "  -3.58%  "
$ws.Range("E10").Value = "  +3.78%  "
$ws.Range("D11").Value = "'5.80"
$ws.Range("E11").Value = "  -14.93%  "
$ws.Range("D12").Value = "'0.385"
$ws.Range("E12").Value = "  -2.95%  "
$ws.Range("D13").Value = "3.253.11"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").Value = "'26.92"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "63.710.12"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("E16").Value = "  -4.87%  "
$ws.Range("D17").Value = "2.766.22"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "'12.14"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").Value = "'4.84"
$ws.Range("E19").Value = "  -3.52%  "
$ws.Range("D20").Value = "'356.15"
$ws.Range("E20").Value = "  -3.27%  "
$ws.Range("E21").Value = "  -5.90%  "
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").Value = "'0.526"
$ws.Range("E23").Value = "  -4.73%  "
$ws.Range("D24").Value = "'65.02"
$ws.Range("E24").Value = "  -3.48%  "
$ws.Range("D25").Value = "'0.169"
$ws.Range("E25").Value = "  -4.23%  "
$ws.Range("D26").Value = "'8.55"
$ws.Range("E26").Value = "  -1.79%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "0.0₃0907"
$ws.Range("E28").Value = "  -6.54%  "
$ws.Range("D29").Value = "'7.37"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -4.73%  "
$ws.Range("D31").Value = "'1.26"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").Value = "'169.81"
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'4.92"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'20.14"
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.47"
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("E37").Value = "  -2.17%  "
$ws.Range("E38").Value = "  -3.18%  "
$ws.Range("D39").Value = "'350.72"
$ws.Range("E39").Value = "  +1.88%  "
$ws.Range("D40").Value = "'6.23"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "'4.17"
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("D42").Value = "'39.09"
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("D43").Value = "'21.48"
$ws.Range("E43").Value = "  -4.63%  "
$ws.Range("D44").Value = "'21.80"
$ws.Range("E44").Value = "  -4.82%  "
$ws.Range("D45").Value = "'0.0587"
$ws.Range("E45").Value = "  -4.36%  "
$ws.Range("D46").Value = "'138.72"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").Value = "'0.0253"
$ws.Range("E48").Value = "  -3.65%  "
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'11.05"
$ws.Range("E51").Value = "  +0.05%  "
